$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The report lists inventory items alphabetically. A new item - a pregnancy
# test kit ("اختبار حمل بيبي تشك") - needs to be inserted as item #17,
# right before "بلاستر مترسيلك 2.5 سم" (which was item #17 and becomes #18),
# pushing every following item down by one row.

# 1) Insert a new blank row at row 20 (pushes old rows 20-25 down to 21-26).
$ws.Rows("20:20").Insert()

# 2) Give the new row the same formatting (styles + merged cells) as the
#    data rows around it by copying row 21 (which now holds the row that
#    used to be row 20) onto the freshly inserted, still-blank row 20.
$ws.Range("A21:N21").Copy()
$ws.Range("A20:N20").PasteSpecial(-4122)
$ws.Range("B20:G20").Merge()
$ws.Range("H20:K20").Merge()
$ws.Range("L20:M20").Merge()

# 3) Populate the new row's data: serial no., name, current balance,
#    sale price and number of transactions.
$ws.Range("A20").Value = 17
$ws.Range("B20").Value = "اختبار حمل بيبي تشك "
$ws.Range("H20").Value = "14:0"
$ws.Range("L20").Value = -25
$ws.Range("N20").Value = "1:0"

# 4) Renumber the "م" (serial number) column for the rows that shifted down.
$ws.Range("A21").Value = 18
$ws.Range("A22").Value = 19
$ws.Range("A23").Value = 20
$ws.Range("A24").Value = 21

# 5) Recompute the displayed total (sum of the sale-price column), which
#    drops by 25 because of the new item's -25 price.
$ws.Range("K25").Value = 1513.04

# 6) Match the final row heights of the shifted block.
$ws.Rows("20:20").RowHeight = 25.5
$ws.Rows("21:21").RowHeight = 24.75
$ws.Rows("22:22").RowHeight = 25.5
$ws.Rows("23:23").RowHeight = 25.5
$ws.Rows("24:24").RowHeight = 24.75
$ws.Rows("25:25").RowHeight = 26.25
$ws.Rows("26:26").RowHeight = 16.5
